$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 1317.7858
$ws.Range("I39").Value = 1181.125
$ws.Range("K39").Value = 3543.375
$ws.Range("M39").Value = -3247.375

$ws.Range("H129").Value = 930.4666999999999
$ws.Range("J129").Value = 915.5797
$ws.Range("L129").Value = 2746.7391
$ws.Range("N129").Value = -12746.7391

$ws.Range("H132").Value = 1121.3077
$ws.Range("I132").Value = 1186.129
$ws.Range("J132").Value = 870.125
$ws.Range("K132").Value = 3558.387
$ws.Range("L132").Value = 2610.375
$ws.Range("M132").Value = -1028.387
$ws.Range("N132").Value = -7670.375

$ws.Range("H135").Value = 556.4
$ws.Range("I135").Value = 502.25
$ws.Range("K135").Value = 4520.25
$ws.Range("M135").Value = -1985.25

$ws.Range("H137").Value = 92642.37
$ws.Range("I137").Value = 800
$ws.Range("J137").Value = 127083.25
$ws.Range("K137").Value = 2400
$ws.Range("L137").Value = 381249.75
$ws.Range("M137").Value = 150
$ws.Range("N137").Value = -386349.75

$ws.Range("H138").Value = 3424.9805
$ws.Range("I138").Value = 4385.067
$ws.Range("J138").Value = 3024.9443
$ws.Range("K138").Value = 13155.201
$ws.Range("L138").Value = 9074.832900000001
$ws.Range("M138").Value = -8015.201000000001
$ws.Range("N138").Value = -19354.8329

$ws.Range("H141").Value = 4355.7144
$ws.Range("I141").Value = 3996.6667
$ws.Range("J141").Value = 4625
$ws.Range("K141").Value = 11990.0001
$ws.Range("L141").Value = 13875
$ws.Range("M141").Value = -6810.000100000001
$ws.Range("N141").Value = -24235

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 932726.75
$ws.Range("I2").Value = 1226167
$ws.Range("K2").Value = 1226167
$ws.Range("M2").Value = -1226054

$ws.Range("H32").Value = 3461.4927
$ws.Range("I32").Value = 2657.5
$ws.Range("J32").Value = 8821.444
$ws.Range("K32").Value = 2657.5
$ws.Range("L32").Value = 8821.444
$ws.Range("M32").Value = -2370.5
$ws.Range("N32").Value = -9395.444

$ws.Range("H102").Value = 1299
$ws.Range("I102").Value = 1299
$ws.Range("K102").Value = 1299
$ws.Range("M102").Value = 323

$ws.Range("H116").Value = 932726.75
$ws.Range("I116").Value = 1226167
$ws.Range("K116").Value = 1226167
$ws.Range("M116").Value = -1223873

$ws.Range("H132").Value = 2444.4482
$ws.Range("I132").Value = 1949.4375
$ws.Range("J132").Value = 3053.6924
$ws.Range("K132").Value = 5848.3125
$ws.Range("L132").Value = 9161.0772
$ws.Range("M132").Value = -3318.3125
$ws.Range("N132").Value = -14221.0772

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 932726.75
$ws.Range("I3").Value = 1226167
$ws.Range("K3").Value = 1226167
$ws.Range("M3").Value = -1226053

$ws.Range("H99").Value = 1110.2
$ws.Range("I99").Value = 1122.5555
$ws.Range("K99").Value = 1122.5555
$ws.Range("M99").Value = 375.4445000000001

$ws.Range("H134").Value = 6789.5454
$ws.Range("I134").Value = 7036.3794
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 21109.1382
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -18574.1382
$ws.Range("N134").Value = -20070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2284.6052
$ws.Range("I31").Value = 1232.4166
$ws.Range("J31").Value = 2770.2307
$ws.Range("K31").Value = 1232.4166
$ws.Range("L31").Value = 2770.2307
$ws.Range("M31").Value = -937.4166
$ws.Range("N31").Value = -3360.2307

$ws.Range("H34").Value = 2284.6052
$ws.Range("I34").Value = 1232.4166
$ws.Range("J34").Value = 2770.2307
$ws.Range("K34").Value = 1232.4166
$ws.Range("L34").Value = 2770.2307
$ws.Range("M34").Value = -1030.4166
$ws.Range("N34").Value = -3174.2307

$ws.Range("H62").Value = 3324.2
$ws.Range("I62").Value = 3290.8572
$ws.Range("K62").Value = 3290.8572
$ws.Range("M62").Value = -2666.8572

$ws.Range("H65").Value = 3324.2
$ws.Range("I65").Value = 3290.8572
$ws.Range("K65").Value = 16454.286
$ws.Range("M65").Value = -13334.286

$ws.Range("H107").Value = 402
$ws.Range("I107").Value = 369
$ws.Range("K107").Value = 369
$ws.Range("M107").Value = 1551

$ws.Range("H122").Value = 1600.7858
$ws.Range("I122").Value = 1119.7
$ws.Range("J122").Value = 2803.5
$ws.Range("K122").Value = 3359.1
$ws.Range("L122").Value = 8410.5
$ws.Range("M122").Value = -909.1000000000004
$ws.Range("N122").Value = -13310.5

$ws.Range("H134").Value = 2433.6956
$ws.Range("I134").Value = 1210.2941
$ws.Range("J134").Value = 5900
$ws.Range("K134").Value = 3630.8823
$ws.Range("L134").Value = 17700
$ws.Range("M134").Value = -1095.8823
$ws.Range("N134").Value = -22770

$ws.Range("H141").Value = 70079.75
$ws.Range("J141").Value = 70079.75
$ws.Range("L141").Value = 70079.75
$ws.Range("N141").Value = -80439.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 295.3871
$ws.Range("I5").Value = 219.875
$ws.Range("K5").Value = 659.625
$ws.Range("M5").Value = -547.625

$ws.Range("H12").Value = 65.46154
$ws.Range("J12").Value = 96.14286
$ws.Range("L12").Value = 288.42858
$ws.Range("N12").Value = -634.42858

$ws.Range("H55").Value = 2348.75
$ws.Range("J55").Value = 2348.75
$ws.Range("L55").Value = 7046.25
$ws.Range("N55").Value = -7400.25

$ws.Range("H68").Value = 2349.2554
$ws.Range("I68").Value = 1152.7333
$ws.Range("J68").Value = 2910.125
$ws.Range("K68").Value = 3458.199900000001
$ws.Range("L68").Value = 8730.375
$ws.Range("M68").Value = -2647.199900000001
$ws.Range("N68").Value = -10352.375

$ws.Range("H71").Value = 2349.2554
$ws.Range("I71").Value = 1152.7333
$ws.Range("J71").Value = 2910.125
$ws.Range("K71").Value = 10374.5997
$ws.Range("L71").Value = 26191.125
$ws.Range("M71").Value = -6318.599700000001
$ws.Range("N71").Value = -34303.125

$ws.Range("H109").Value = 4470.143
$ws.Range("I109").Value = 1618
$ws.Range("J109").Value = 6054.6665
$ws.Range("K109").Value = 4854
$ws.Range("L109").Value = 18163.9995
$ws.Range("M109").Value = -3814
$ws.Range("N109").Value = -20243.9995

$ws.Range("H131").Value = 9629679
$ws.Range("J131").Value = 15760.639
$ws.Range("L131").Value = 47281.917
$ws.Range("N131").Value = -57361.917

$ws.Range("H132").Value = 1724.75

$ws.Range("H135").Value = 295.3871
$ws.Range("I135").Value = 219.875
$ws.Range("K135").Value = 1978.875
$ws.Range("M135").Value = 556.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 17666.666

$ws.Range("H73").Value = 17666.666

$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()

$ws.Range("H102").Value = 2401
$ws.Range("I102").Value = 2501.2222
$ws.Range("K102").Value = 2501.2222
$ws.Range("M102").Value = -879.2222000000002

$ws.Range("H110").Value = 99741
$ws.Range("J110").Value = 99741
$ws.Range("L110").Value = 99741
$ws.Range("N110").Value = -107921

$ws.Range("H132").Value = 2027816.4
$ws.Range("I132").Value = 3207372.5
$ws.Range("J132").Value = 5720.2856
$ws.Range("K132").Value = 9622117.5
$ws.Range("L132").Value = 17160.8568
$ws.Range("M132").Value = -9619587.5
$ws.Range("N132").Value = -22220.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2115.2856
$ws.Range("I22").Value = 1800
$ws.Range("J22").Value = 2167.8333
$ws.Range("K22").Value = 1800
$ws.Range("L22").Value = 2167.8333
$ws.Range("M22").Value = -1505
$ws.Range("N22").Value = -2757.8333

$ws.Range("H27").Value = 2115.2856
$ws.Range("I27").Value = 1800
$ws.Range("J27").Value = 2167.8333
$ws.Range("K27").Value = 1800
$ws.Range("L27").Value = 2167.8333
$ws.Range("M27").Value = -1693
$ws.Range("N27").Value = -2381.8333

$ws.Range("H46").Value = 1055.4445
$ws.Range("I46").Value = 500.83334
$ws.Range("J46").Value = 1332.75
$ws.Range("K46").Value = 500.83334
$ws.Range("L46").Value = 1332.75
$ws.Range("M46").Value = -312.83334
$ws.Range("N46").Value = -1708.75

$ws.Range("H136").Value = 2707.3057
$ws.Range("I136").Value = 1375.72
$ws.Range("J136").Value = 5733.636
$ws.Range("K136").Value = 4127.16
$ws.Range("L136").Value = 17200.908
$ws.Range("M136").Value = -1577.16
$ws.Range("N136").Value = -22300.908

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws.Range("H123").Value = 47408.535
$ws.Range("J123").Value = 47408.535
$ws.Range("L123").Value = 47408.535
$ws.Range("N123").Value = -57208.535

$ws.Range("H136").Value = 8171918
$ws.Range("I136").Value = 12078933
$ws.Range("J136").Value = 2704.4092
$ws.Range("K136").Value = 36236799
$ws.Range("L136").Value = 8113.2276
$ws.Range("M136").Value = -36234249
$ws.Range("N136").Value = -13213.2276
